# Updated queries for C3DC first half testcases.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the SQL JOIN/ON clauses in every query cell on the sheet.
#    The "id" / "study.id" / "participant.id" columns were renamed to the
#    fully-qualified "study_id" / "study.study_id" / "participant_id" /
#    "participant.participant_id" column names used by the new data model.
# ---------------------------------------------------------------------------
$queryCells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $queryCells) {
    $rng = $ws.Range($addr)
    $text = $rng.Text

    $text = $text.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $text = $text.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')

    $rng.Value = $text
}

# ---------------------------------------------------------------------------
# 2) Scroll the sheet view so column A (not B) is the left-most visible
#    column, keeping the current selection (C7) untouched.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1

# ---------------------------------------------------------------------------
# 3) Widen column C (drop its old "best fit" width) to fit the longer query
#    text, and grow row 2's height to match the rewrapped content.
# ---------------------------------------------------------------------------
$ws.Range("C1").ColumnWidth = 67.5
$ws.Range("A2").RowHeight = 336
